$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.886.91"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.623.74"
$ws.Range("E3").Value = "  -0.37%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.42"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.515"
$ws.Range("E6").Value = "  -1.58%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.97"
$ws.Range("E8").Value = "  -1.51%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.257"
$ws.Range("E9").Value = "  +0.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0604"
$ws.Range("E10").Value = "  -1.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0880"
$ws.Range("E11").Value = "  -0.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.857.57"
$ws.Range("E12").Value = "  -0.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.640.16"
$ws.Range("E13").Value = "  +0.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.00"
$ws.Range("E14").Value = "  -0.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.552"
$ws.Range("E15").Value = "  -1.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.52"
$ws.Range("E16").Value = "  -1.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.903.43"
$ws.Range("E17").Value = "  -0.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.37"
$ws.Range("E18").Value = "  -1.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.59"
$ws.Range("E19").Value = "  -0.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0715"
$ws.Range("E20").Value = "  -1.18%  "

$ws.Range("E21").Value = "  +0.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.33"
$ws.Range("E22").Value = "  -0.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.93"
$ws.Range("E23").Value = "  -3.31%  "

$ws.Range("E24").Value = "  +1.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.31"
$ws.Range("E25").Value = "  -0.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.90"
$ws.Range("E26").Value = "  -0.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.110"
$ws.Range("E27").Value = "  -0.67%  "

$ws.Range("E28").Value = "  +0.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.34"
$ws.Range("E29").Value = "  -1.34%  "

$ws.Range("E30").Value = "  -0.33%  "

$ws.Range("E31").Value = "  -0.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.39"
$ws.Range("E32").Value = "  -0.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.412.05"
$ws.Range("E33").Value = "  +0.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.09"
$ws.Range("E34").Value = "  +0.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("E35").Value = "  +1.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.996"
$ws.Range("E36").Value = "  -2.20%  "

$ws.Range("E37").Value = "  -0.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0169"
$ws.Range("E38").Value = "  -1.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.552"
$ws.Range("E39").Value = "  -0.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.847"
$ws.Range("E40").Value = "  -1.91%  "

$ws.Range("E41").Value = "  +0.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -2.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.81"
$ws.Range("E43").Value = "  -1.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.21"
$ws.Range("E44").Value = "  -1.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.38"
$ws.Range("E45").Value = "  -2.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.765.35"
$ws.Range("E46").Value = "  -0.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.83"
$ws.Range("E48").Value = "  +0.89%  "

$ws.Range("E49").Value = "  -0.73%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.100"
$ws.Range("E50").Value = "  -0.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0502"
$ws.Range("E51").Value = "  -0.55%  "
